$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Regression-coefficient table (rows 17-23) ------------------------------
# The old sheet baked the quadratic-fit coefficients directly into the F17:F23
# formulas ( -0.0002278815 / (-249.575444 + 0.6036185*C + -0.000365585*C*C) ).
# They are pulled out into labelled cells F17:F20 = "a","b","c","d" with the
# actual coefficient values placed in G17:G20, and H17 now estimates the
# scaling factor at a y-distance of 2297 px using those named coefficients.

$ws.Range("F17").Value = "a"
$ws.Range("F18").Value = "b"
$ws.Range("F19").Value = "c"
$ws.Range("F20").Value = "d"

# Labels are right aligned (they previously held right-aligned numbers).
$ws.Range("F17:F20").HorizontalAlignment = -4152   # xlRight

$ws.Range("G17").Value = -0.00022788150560381401
$ws.Range("G18").Value = -249.575444271701
$ws.Range("G19").Value = 0.60361848274124597
$ws.Range("G20").Value = -0.00036558500089469299

$ws.Range("H17").Formula = "=G17/(G18+G19*(2297)+G20*(2297*2297))"

# The old per-row "fit vs. measured" comparison (F18:F23 / G18:G23) is gone;
# only F21:F23 remain as blank, still-numeric-formatted cells.
$ws.Range("F21:F23").ClearContents()
$ws.Range("G21:G23").ClearContents()

# --- Selection ---------------------------------------------------------------
$ws.Range("D17").Select() | Out-Null
